$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simulation")

$ws.Range("B2").Value = 0.31096858330410443
$ws.Range("C2").Value = -0.8654311929324662
$ws.Range("D2").Value = 0.15694806387753057
$ws.Range("E2").Value = 0.08537898030687405
